$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (columns A:R)
# A new salesperson/sub-agent "MULLO GUACHO ANA LUCIA" is inserted as row 14
# for asesor "ILLER LOPEZ ROBERTO FERNANDO", pushing the existing rows for
# "PAUTA ASTUDILLO JULIO HERNAN" and "VIEJO RIVAS MAYRA ANABELLE" down by one
# row, and the trailing totals row's "X de 14" labels become "X de 15".
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows.Item(14).Insert()

$ws1.Cells.Item(14, 1).Value() = "ILLER LOPEZ ROBERTO FERNANDO"
$ws1.Cells.Item(14, 2).Value() = "MULLO GUACHO ANA LUCIA"
for ($col = 3; $col -le 18; $col++) {
  $ws1.Cells.Item(14, $col).Value() = 0
}

for ($col = 3; $col -le 18; $col++) {
  $cell = $ws1.Cells.Item(17, $col)
  $cell.Value() = $cell.Value().Replace("de 14", "de 15")
}

# ---------------------------------------------------------------------------
# Sheet "VENTA MENSUAL" (columns A:G) gets the same row inserted.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(14).Insert()

$ws2.Cells.Item(14, 1).Value() = "ILLER LOPEZ ROBERTO FERNANDO"
$ws2.Cells.Item(14, 2).Value() = "MULLO GUACHO ANA LUCIA"
for ($col = 3; $col -le 7; $col++) {
  $ws2.Cells.Item(14, $col).Value() = 0
}
